$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.782.25'
$ws.Range("E2").Value = '  -0.24%  '

# Row 3
$ws.Range("D3").Value = '2.290.04'
$ws.Range("E3").Value = '  -0.21%  '

# Row 4
$ws.Range("E4").Value = '  -0.24%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '115.09'
$ws.Range("E5").Value = '  +2.28%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '266.77'
$ws.Range("E6").Value = '  -1.12%  '

# Row 7
$ws.Range("E7").Value = '  +3.72%  '

# Row 8
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  -0.68%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.55'
$ws.Range("E10").Value = '  -0.90%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0941'
$ws.Range("E11").Value = '  -0.79%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.22'
$ws.Range("E12").Value = '  +0.99%  '

# Row 13
$ws.Range("E13").Value = '  +1.42%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.46'
$ws.Range("E14").Value = '  -2.56%  '

# Row 15
$ws.Range("D15").Value = '2.631.39'
$ws.Range("E15").Value = '  -0.29%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.880'
$ws.Range("E16").Value = '  +3.19%  '

# Row 17
$ws.Range("D17").Value = '2.285.95'
$ws.Range("E17").Value = '  -0.62%  '

# Row 18
$ws.Range("D18").Value = '43.633.32'
$ws.Range("E18").Value = '  -0.25%  '

# Row 19
$ws.Range("E19").Value = '  +0.52%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.90'
$ws.Range("E20").Value = '  +1.94%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.55'
$ws.Range("E21").Value = '  +0.32%  '

# Row 22
$ws.Range("E22").Value = '  +0.01%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.73'
$ws.Range("E23").Value = '  +1.72%  '

# Row 24
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.57'
$ws.Range("E24").Value = '  -1.91%  '

# Row 25
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.90'
$ws.Range("E25").Value = '  +2.28%  '

# Row 26
$ws.Range("E26").Value = '  +1.74%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.64'
$ws.Range("E27").Value = '  -0.42%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.93'
$ws.Range("E28").Value = '  -0.05%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.42'
$ws.Range("E29").Value = '  +0.66%  '

# Row 30
$ws.Range("E30").Value = '  -0.85%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.67'
$ws.Range("E31").Value = '  -1.30%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.77'
$ws.Range("E32").Value = '  +0.94%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0917'
$ws.Range("E33").Value = '  -1.94%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.72'
$ws.Range("E34").Value = '  +0.65%  '

# Row 35
$ws.Range("E35").Value = '  +1.87%  '

# Row 36
$ws.Range("E36").Value = '  +5.63%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.70'
$ws.Range("E37").Value = '  +0.79%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.96'
$ws.Range("E38").Value = '  +3.44%  '

# Row 39
$ws.Range("E39").Value = '  -0.93%  '

# Row 40
$ws.Range("E40").Value = '  +6.60%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.28'
$ws.Range("E41").Value = '  +4.41%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.45'
$ws.Range("E42").Value = '  -0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.236'
$ws.Range("E43").Value = '  -2.97%  '

# Row 44
$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.99'
$ws.Range("E44").Value = '  -6.16%  '

# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.12%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.38'
$ws.Range("E46").Value = '  -0.91%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '75.08'
$ws.Range("E47").Value = '  +36.08%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.30'
$ws.Range("E48").Value = '  +4.64%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.62'
$ws.Range("E49").Value = '  -1.79%  '

# Row 50
$ws.Range("E50").Value = '  +0.49%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '101.12'
$ws.Range("E51").Value = '  -1.89%  '
